# Generate Report for Handback
# Updates row 7 (the 8f0df88c-7249-4638-947e-33f48217245b file) on both the
# "zh-cn" and "de-de" sheets: the handback has now been picked up, but it is
# not the latest version, so:
#   - I7 (Latest Target File)      gets a hyperlink to the handback .md file
#   - J7 (Latest Handback File)    gets the generated .xlf file name
#   - K7 (Latest Handback DateTime) gets the generation timestamp
#   - P7 (Error Detail)            gets the "not latest version" message

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/405f5638f604a5f64651aaf652078c57fcf91e6c/e2e/8f0df88c-7249-4638-947e-33f48217245b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc4d47ec8b778b3a0b6ed7c2227a67c689a741e2/e2e/8f0df88c-7249-4638-947e-33f48217245b.md."

$handbackDisplay = "8f0df88c-7249-4638-947e-33f48217245b.md"

function Update-HandbackRow {
    param($ws, $langSlug, $handoffXlf, $handbackDateTime)

    # Latest Target File: add the hyperlink to the handback markdown file.
    $ws.Hyperlinks.Add(
        $ws.Range("I7"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0-$langSlug/blob/405f5638f604a5f64651aaf652078c57fcf91e6c/e2e/8f0df88c-7249-4638-947e-33f48217245b.md",
        "",
        "",
        $handbackDisplay
    )

    # Latest Handback File
    $ws.Range("J7").Value = $handoffXlf

    # Latest Handback DateTime
    $ws.Range("K7").Value = $handbackDateTime

    # Error Detail
    $ws.Range("P7").Value = $errorDetail
}

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "zhcn" "8f0df88c-7249-4638-947e-33f48217245b.82067bd792b9b870f82f6d61471e19e81e7aa19c.zh-cn.xlf" "2016-08-24 02:54:52"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "dede" "8f0df88c-7249-4638-947e-33f48217245b.82067bd792b9b870f82f6d61471e19e81e7aa19c.de-de.xlf" "2016-08-24 02:54:59"
